$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)
$s.Delete()
